$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 45033
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 24000
$ws.Cells.Item(2, 12).Value = 24000
$ws.Cells.Item(2, 13).Value = 24000
$ws.Cells.Item(2, 14).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(2, 16).Value = 1333

# Row 3
$ws.Cells.Item(3, 4).Value = 45037
$ws.Cells.Item(3, 10).Value = 80
$ws.Cells.Item(3, 11).Value = 24000
$ws.Cells.Item(3, 12).Value = 24000
$ws.Cells.Item(3, 13).Value = 24000
$ws.Cells.Item(3, 14).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(3, 16).Value = 1600
$ws.Cells.Item(3, 17).Value = 15

# Row 4
$ws.Cells.Item(4, 4).Value = 45014
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 8000
$ws.Cells.Item(4, 13).Value = 8000
$ws.Cells.Item(4, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(4, 16).Value = 444
$ws.Cells.Item(4, 17).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 45034
$ws.Cells.Item(5, 10).Value = 50
$ws.Cells.Item(5, 11).Value = 24000
$ws.Cells.Item(5, 12).Value = 24000
$ws.Cells.Item(5, 13).Value = 24000
$ws.Cells.Item(5, 14).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(5, 16).Value = 1333

# Row 6
$ws.Cells.Item(6, 4).Value = 45015
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(6, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(6, 16).Value = 1333
$ws.Cells.Item(6, 17).Value = 18

# Row 7
$ws.Cells.Item(7, 4).Value = 44280
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 30
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 25000
$ws.Cells.Item(7, 13).Value = 25000
$ws.Cells.Item(7, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(7, 16).Value = 1389
$ws.Cells.Item(7, 17).Value = 18

# Row 8
$ws.Cells.Item(8, 4).Value = 44313
$ws.Cells.Item(8, 14).Value = '$/caja 15 kilos empedrada'

# Row 9
$ws.Cells.Item(9, 4).Value = 44313
$ws.Cells.Item(9, 10).Value = 20
$ws.Cells.Item(9, 11).Value = 30000
$ws.Cells.Item(9, 12).Value = 30000
$ws.Cells.Item(9, 13).Value = 30000
$ws.Cells.Item(9, 14).Value = '$/caja 20 kilos empedrada'
$ws.Cells.Item(9, 16).Value = 1500
$ws.Cells.Item(9, 17).Value = 20

# Row 10
$ws.Cells.Item(10, 4).Value = 45042
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 24000
$ws.Cells.Item(10, 13).Value = 24000
$ws.Cells.Item(10, 16).Value = 1333

# Row 11
$ws.Cells.Item(11, 4).Value = 45041
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 24000
$ws.Cells.Item(11, 13).Value = 24000
$ws.Cells.Item(11, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(11, 16).Value = 1333
$ws.Cells.Item(11, 17).Value = 18

# Row 12
$ws.Cells.Item(12, 4).Value = 44285
$ws.Cells.Item(12, 11).Value = 25000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 25000
$ws.Cells.Item(12, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(12, 16).Value = 1389
$ws.Cells.Item(12, 17).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = 45040
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 18000
$ws.Cells.Item(13, 14).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(13, 16).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = 44315
$ws.Cells.Item(14, 9).Value = 'Especial'
$ws.Cells.Item(14, 10).Value = 10
$ws.Cells.Item(14, 11).Value = 30000
$ws.Cells.Item(14, 12).Value = 30000
$ws.Cells.Item(14, 13).Value = 30000
$ws.Cells.Item(14, 14).Value = '$/caja 20 kilos empedrada'
$ws.Cells.Item(14, 16).Value = 1500
$ws.Cells.Item(14, 17).Value = 20

# Row 15
$ws.Cells.Item(15, 4).Value = 44315
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 15000
$ws.Cells.Item(15, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(15, 16).Value = 1000
$ws.Cells.Item(15, 17).Value = 15

# New row 16 (full new record)
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(16, 3).Value = 'La Araucanía'
$ws.Cells.Item(16, 4).Value = 44293
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100112041
$ws.Cells.Item(16, 7).Value = 'Fruto del paraíso'
$ws.Cells.Item(16, 8).Value = 'Sin especificar'
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 10
$ws.Cells.Item(16, 11).Value = 25000
$ws.Cells.Item(16, 12).Value = 25000
$ws.Cells.Item(16, 13).Value = 25000
$ws.Cells.Item(16, 14).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(16, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(16, 16).Value = 1667
$ws.Cells.Item(16, 17).Value = 15
$ws.Cells.Item(16, 18).Value = 'Hortaliza'

# Apply the date style/number format used by column D to the new row
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat